{"js": "// The protocol's JSON example snippets describe each \"item\" object as\n// {'id' : <num>, 'amount' : <num>}. This edit documents that 'id' is\n// actually a string, not a number, so every \"<num>, \" that immediately\n// follows \"'id' :\" becomes \"<string>, \" \u2014 the sibling \"'amount' : <num>\"\n// stays a number and is left untouched.\nconst body = context.document.body;\n\n// \"<num>, \" (angle bracket + num + angle bracket + comma + space) only ever\n// occurs right after the 'id' key in these snippets (the 'amount' key is\n// always the last field in the object, so it is followed by \"}\" instead of\n// \", \"). Searching for this exact substring therefore targets only the\n// six 'id' field-type declarations that need to change.\nconst matches = body.search(\"<num>, \", { matchCase: true });\nmatches.load(\"text\");\nawait context.sync();\n\nfor (const match of matches.items) {\n  match.insertText(\"<string>, \", \"Replace\");\n}\nawait context.sync();\n", "ps1": "# The protocol's JSON example snippets describe each \"item\" object as\n# {'id' : <num>, 'amount' : <num>}. This edit documents that 'id' is\n# actually a string, not a number, so every \"<num>, \" that immediately\n# follows 'id' becomes \"<string>, \" \u2014 the sibling 'amount' : <num> stays\n# a number and is left untouched.\n#\n# \"<num>, \" (angle bracket + num + angle bracket + comma + space) only\n# ever occurs right after the 'id' key in these snippets (the 'amount'\n# key is always the last field in the object, so it is followed by \"}\"\n# instead of \", \"). A scoped Find/Replace on that exact substring\n# therefore only touches the six 'id' field-type declarations, leaving\n# every 'amount' : <num> untouched.\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Replacement.ClearFormatting()\n\n# Find.Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,\n#   MatchSoundsLike, MatchAllWordForms, Forward, Wrap, Format,\n#   ReplaceWith, Replace)\n$d.Content.Find.Execute(\n    \"<num>, \",      # FindText\n    $true,          # MatchCase\n    $false,         # MatchWholeWord\n    $false,         # MatchWildcards\n    $false,         # MatchSoundsLike\n    $false,         # MatchAllWordForms\n    $true,          # Forward\n    1,              # Wrap -> wdFindContinue\n    $false,         # Format\n    \"<string>, \",   # ReplaceWith\n    2               # Replace -> wdReplaceAll\n) | Out-Null\n"}
